$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D6 and D8 get new "Done" status values (Create / Read tests for Accounts module)
$ws.Range("D6").Value = "Done"
$ws.Range("D8").Value = "Done"

# Highlight the Status cells for the Accounts "Create" and "Read" rows green
$ws.Range("D4").Interior.Color = 5296274
$ws.Range("D5").Interior.Color = 5296274
$ws.Range("D6").Interior.Color = 5296274
$ws.Range("D8").Interior.Color = 5296274

# Move the active selection to D9
$ws.Range("D9").Select() | Out-Null
